$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-11 Wednesday" "2023-10-12 Thursday"

Replace-Text "31÷2=" "53÷2="
Replace-Text "89÷8=" "96÷7="
Replace-Text "87÷6=" "22÷3="
Replace-Text "78÷8=" "32÷8="
Replace-Text "84÷9=" "77÷5="

Replace-Text "57÷2=" "41÷8="
Replace-Text "48÷5=" "49÷7="
Replace-Text "13÷9=" "37÷5="
Replace-Text "52÷4=" "41÷6="
Replace-Text "42÷8=" "71÷8="

Replace-Text "53÷5=" "56÷8="
Replace-Text "69÷2=" "95÷5="
Replace-Text "63÷6=" "47÷5="
Replace-Text "21÷6=" "22÷2="
Replace-Text "45÷6=" "11÷4="

Replace-Text "19÷4=" "23÷3="
Replace-Text "51÷2=" "43÷6="
Replace-Text "53÷8=" "89÷8="
Replace-Text "16÷9=" "41÷7="
Replace-Text "57÷5=" "26÷8="

Replace-Text "73÷3=" "85÷5="
Replace-Text "90÷6=" "21÷5="
Replace-Text "41÷2=" "32÷3="
Replace-Text "28÷3=" "39÷3="
Replace-Text "24÷4=" "16÷4="
